# Auto-generated edit script applying 88 cell changes to Sheet1
# Source: GitHub Actions crypto price refresh commit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $CellRef, $NewValue)
    # Force the value to be stored as text (matches the source data's
    # inline-string cell type) by using a leading apostrophe, then clear
    # the resulting quote-prefix style so formatting stays untouched.
    $Sheet.Range($CellRef).Value = "'" + $NewValue
    $Sheet.Range($CellRef).Style = "Normal"
}

Set-TextCell $ws 'D2' '59.047.51'
$ws.Range("E2").Value = '  -4.75%  '
Set-TextCell $ws 'D3' '2.517.10'
$ws.Range("E3").Value = '  -2.43%  '
Set-TextCell $ws 'D4' '0.998'
$ws.Range("E4").Value = '  -0.22%  '
Set-TextCell $ws 'D5' '536.12'
$ws.Range("E5").Value = '  -2.36%  '
Set-TextCell $ws 'D6' '145.34'
$ws.Range("E6").Value = '  -6.02%  '
Set-TextCell $ws 'D7' '0.997'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("E8").Value = '  -3.23%  '
Set-TextCell $ws 'D9' '2.515.91'
$ws.Range("E9").Value = '  -2.65%  '
Set-TextCell $ws 'D10' '0.0996'
$ws.Range("E10").Value = '  -4.56%  '
$ws.Range("E11").Value = '  -2.65%  '
Set-TextCell $ws 'D12' '5.58'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("E13").Value = '  -3.34%  '
Set-TextCell $ws 'D14' '2.926.35'
$ws.Range("E14").Value = '  -3.65%  '
Set-TextCell $ws 'D15' '24.01'
$ws.Range("E15").Value = '  -6.16%  '
Set-TextCell $ws 'D16' '58.974.81'
$ws.Range("E16").Value = '  -4.80%  '
$ws.Range("E17").Value = '  -3.73%  '
Set-TextCell $ws 'D18' '2.502.57'
$ws.Range("E18").Value = '  -3.08%  '
Set-TextCell $ws 'D19' '11.29'
$ws.Range("E19").Value = '  -3.21%  '
$ws.Range("E20").Value = '  -5.52%  '
Set-TextCell $ws 'D21' '323.84'
$ws.Range("E21").Value = '  -4.19%  '
Set-TextCell $ws 'D22' '0.998'
$ws.Range("E22").Value = '  -0.11%  '
Set-TextCell $ws 'D23' '5.75'
$ws.Range("E23").Value = '  -4.77%  '
Set-TextCell $ws 'D24' '61.33'
$ws.Range("E24").Value = '  -3.79%  '
Set-TextCell $ws 'D25' '0.441'
$ws.Range("E25").Value = '  -10.60%  '
$ws.Range("E26").Value = '  -3.82%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D27' '0.994'
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell $ws 'D28' '2.605.83'
$ws.Range("E28").Value = '  -3.50%  '
Set-TextCell $ws 'D29' '7.79'
$ws.Range("E29").Value = '  -4.68%  '
Set-TextCell $ws 'D30' '6.88'
$ws.Range("E30").Value = '  -5.69%  '
$ws.Range("E31").Value = '  -6.88%  '
$ws.Range("E32").Value = '  -6.65%  '
$ws.Range("E33").Value = '  -5.36%  '
$ws.Range("E34").Value = '  -0.27%  '
Set-TextCell $ws 'D35' '158.60'
$ws.Range("E35").Value = '  -2.67%  '
$ws.Range("E36").Value = '  +2.15%  '
Set-TextCell $ws 'D37' '18.53'
$ws.Range("E37").Value = '  -3.52%  '
$ws.Range("E38").Value = '  -8.84%  '
$ws.Range("E39").Value = '  -9.77%  '
Set-TextCell $ws 'D40' '5.92'
$ws.Range("E40").Value = '  -1.90%  '
Set-TextCell $ws 'D41' '307.90'
$ws.Range("E41").Value = '  -6.48%  '
Set-TextCell $ws 'D42' '36.79'
$ws.Range("E42").Value = '  -2.25%  '
Set-TextCell $ws 'D43' '3.69'
$ws.Range("E43").Value = '  -6.50%  '
Set-TextCell $ws 'D44' '0.825'
$ws.Range("E44").Value = '  -9.22%  '
Set-TextCell $ws 'D45' '0.998'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("E47").Value = '  -1.45%  '
Set-TextCell $ws 'D48' '124.31'
$ws.Range("E48").Value = '  +0.91%  '
Set-TextCell $ws 'D49' '0.0931'
$ws.Range("E49").Value = '  -3.48%  '
Set-TextCell $ws 'D50' '18.62'
$ws.Range("E50").Value = '  -4.71%  '
$ws.Range("E51").Value = '  -5.52%  '
